$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.144578313253012
$ws.Range("C2").Value = 0.6385542168674698
$ws.Range("J2").Value = 0.01606425702811245
$ws.Range("O2").Value = 0.004016064257028112
$ws.Range("P2").Value = 0.1124497991967871
$ws.Range("S2").Value = 0.08433734939759036
$ws.Range("C3").Value = 0.05789473684210526
$ws.Range("J3").Value = 0.02105263157894737
$ws.Range("P3").Value = 0.8052631578947368
$ws.Range("S3").Value = 0.1157894736842105
$ws.Range("J4").Value = 0.1132075471698113
$ws.Range("P4").Value = 0.660377358490566
$ws.Range("S4").Value = 0.2264150943396226
$ws.Range("B6").Value = 0.04301075268817205
$ws.Range("D6").Value = 0.02150537634408602
$ws.Range("F6").Value = 0.05376344086021505
$ws.Range("J6").Value = 0.2526881720430108
$ws.Range("O6").Value = 0.01075268817204301
$ws.Range("Q6").Value = 0.2311827956989247
$ws.Range("R6").Value = 0.07526881720430108
$ws.Range("S6").Value = 0.3118279569892473
$ws.Range("B7").Value = 0.09550561797752809
$ws.Range("D7").Value = 0.01685393258426966
$ws.Range("F7").Value = 0.05617977528089887
$ws.Range("J7").Value = 0.151685393258427
$ws.Range("O7").Value = 0.01685393258426966
$ws.Range("Q7").Value = 0.1741573033707865
$ws.Range("R7").Value = 0.06741573033707865
$ws.Range("S7").Value = 0.4213483146067415
$ws.Range("B8").Value = 0.08974358974358974
$ws.Range("D8").Value = 0.03076923076923077
$ws.Range("E8").Value = 0.002564102564102564
$ws.Range("F8").Value = 0.05641025641025641
$ws.Range("J8").Value = 0.1025641025641026
$ws.Range("O8").Value = 0.01794871794871795
$ws.Range("Q8").Value = 0.2051282051282051
$ws.Range("R8").Value = 0.1307692307692308
$ws.Range("S8").Value = 0.3641025641025641
$ws.Range("B9").Value = 0.07425742574257425
$ws.Range("D9").Value = 0.009900990099009901
$ws.Range("F9").Value = 0.04950495049504951
$ws.Range("J9").Value = 0.1188118811881188
$ws.Range("O9").Value = 0.0198019801980198
$ws.Range("Q9").Value = 0.2227722772277228
$ws.Range("R9").Value = 0.09900990099009901
$ws.Range("S9").Value = 0.405940594059406
$ws.Range("B10").Value = 0.1057692307692308
$ws.Range("D10").Value = 0.02483974358974359
$ws.Range("F10").Value = 0.05849358974358974
$ws.Range("J10").Value = 0.125
$ws.Range("O10").Value = 0.01442307692307692
$ws.Range("Q10").Value = 0.217948717948718
$ws.Range("R10").Value = 0.09615384615384616
$ws.Range("S10").Value = 0.3573717948717949
$ws.Range("G11").Value = 0.1277372262773723
$ws.Range("J11").Value = 0.1021897810218978
$ws.Range("K11").Value = 0.1715328467153285
$ws.Range("L11").Value = 0.5912408759124088
$ws.Range("S11").Value = 0.0072992700729927
$ws.Range("G12").Value = 0.7239263803680982
$ws.Range("J12").Value = 0.2392638036809816
$ws.Range("L12").Value = 0.01840490797546012
$ws.Range("S12").Value = 0.01840490797546012
$ws.Range("G13").Value = 0.7428571428571429
$ws.Range("J13").Value = 0.2571428571428571
$ws.Range("F15").Value = 0.03191489361702127
$ws.Range("H15").Value = 0.1542553191489362
$ws.Range("I15").Value = 0.07446808510638298
$ws.Range("J15").Value = 0.324468085106383
$ws.Range("K15").Value = 0.09042553191489362
$ws.Range("M15").Value = 0.005319148936170213
$ws.Range("O15").Value = 0.05851063829787234
$ws.Range("S15").Value = 0.2606382978723404
$ws.Range("F16").Value = 0.004854368932038835
$ws.Range("H16").Value = 0.1407766990291262
$ws.Range("I16").Value = 0.05825242718446602
$ws.Range("J16").Value = 0.4077669902912621
$ws.Range("K16").Value = 0.1650485436893204
$ws.Range("M16").Value = 0.02912621359223301
$ws.Range("N16").Value = 0.004854368932038835
$ws.Range("O16").Value = 0.05825242718446602
$ws.Range("S16").Value = 0.1310679611650485
$ws.Range("F17").Value = 0.01252609603340292
$ws.Range("H17").Value = 0.1920668058455115
$ws.Range("I17").Value = 0.08768267223382047
$ws.Range("J17").Value = 0.4237995824634656
$ws.Range("K17").Value = 0.1148225469728601
$ws.Range("M17").Value = 0.01461377870563674
$ws.Range("O17").Value = 0.05636743215031315
$ws.Range("S17").Value = 0.09812108559498957
$ws.Range("F18").Value = 0.01834862385321101
$ws.Range("H18").Value = 0.1513761467889908
$ws.Range("I18").Value = 0.07798165137614679
$ws.Range("J18").Value = 0.5091743119266054
$ws.Range("K18").Value = 0.09174311926605505
$ws.Range("M18").Value = 0.01376146788990826
$ws.Range("O18").Value = 0.04128440366972477
$ws.Range("S18").Value = 0.0963302752293578
$ws.Range("F19").Value = 0.01702508960573477
$ws.Range("H19").Value = 0.1872759856630824
$ws.Range("I19").Value = 0.1057347670250896
$ws.Range("J19").Value = 0.3933691756272402
$ws.Range("K19").Value = 0.09408602150537634
$ws.Range("M19").Value = 0.01792114695340502
$ws.Range("N19").Value = 0.0008960573476702509
$ws.Range("O19").Value = 0.06720430107526881
$ws.Range("S19").Value = 0.1164874551971326

Write-Output "Applied team specific time data updates"
